# Revert "I added a mistake on purpose"
#
# The document had an extra paragraph ("This is clearly another
# mistake!") appended at the end, and the `_GoBack` bookmark that used
# to sit at the end of the "... Mouse PFKAM_MOUSE P47857" paragraph had
# been moved onto that extra paragraph. Reverting means: move the
# `_GoBack` bookmark back to the start of the "Mouse PFKAM_MOUSE
# P47857" paragraph, and delete the extra "This is clearly another
# mistake!" paragraph entirely (including its paragraph mark).

$d = $word.ActiveDocument

# Locate the paragraph that ends with the PFKAM_MOUSE text and the
# paragraph that holds the bogus sentence.
$pfkmMouse = $null
$mistake = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text -like "*PFKAM_MOUSE P47857*") {
        $pfkmMouse = $para
    }
    if ($text -like "*This is clearly another mistake!*") {
        $mistake = $para
    }
}

# Re-create the `_GoBack` bookmark at the very start of the
# PFKM mouse paragraph (this is where it lived before the mistake was
# introduced).
$bmRange = $pfkmMouse.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the whole erroneous paragraph, including its paragraph mark.
$mistake.Range.Delete()
